# Update column F (dSF) values for rows 4-9 on Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F4").Value = 9
$ws.Range("F5").Value = -7
$ws.Range("F6").Value = 0
$ws.Range("F7").Value = -2
$ws.Range("F8").Value = -1
$ws.Range("F9").Value = -2
